$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data to latest scrape

# Row 2
$ws.Cells.Item(2, 4).Value = "31.473.20"
$ws.Cells.Item(2, 5).Value = "  +3.79%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "1.990.84"
$ws.Cells.Item(3, 5).Value = "  +6.36%  "

# Row 4
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "0.9993"
$ws.Cells.Item(4, 5).Value = "  -0.08%  "

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "0.8181"
$ws.Cells.Item(5, 5).Value = "  +74.23%  "

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "253.90"
$ws.Cells.Item(6, 5).Value = "  +4.22%  "

# Row 7
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.9993"
$ws.Cells.Item(7, 5).Value = "  -0.08%  "

# Row 8
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.3417"
$ws.Cells.Item(8, 5).Value = "  +18.82%  "

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "25.80"
$ws.Cells.Item(9, 5).Value = "  +16.69%  "

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "0.07026"
$ws.Cells.Item(10, 5).Value = "  +9.02%  "

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.8481"
$ws.Cells.Item(11, 5).Value = "  +17.72%  "

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "0.08122"
$ws.Cells.Item(12, 5).Value = "  +4.56%  "

# Row 13
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "102.16"
$ws.Cells.Item(13, 5).Value = "  +6.91%  "

# Row 14
$ws.Cells.Item(14, 4).Value = "1.987.54"
$ws.Cells.Item(14, 5).Value = "  +6.20%  "

# Row 15
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "5.497"
$ws.Cells.Item(15, 5).Value = "  +7.22%  "

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "275.43"
$ws.Cells.Item(16, 5).Value = "  -1.45%  "

# Row 17
$ws.Cells.Item(17, 4).Value = "31.468.05"
$ws.Cells.Item(17, 5).Value = "  +3.81%  "

# Row 18
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "14.01"
$ws.Cells.Item(18, 5).Value = "  +7.85%  "

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "0.000007979"
$ws.Cells.Item(19, 5).Value = "  +7.27%  "

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "5.730"
$ws.Cells.Item(20, 5).Value = "  +9.55%  "

# Row 21
$ws.Cells.Item(21, 4).Value = "2.252.83"
$ws.Cells.Item(21, 5).Value = "  +6.28%  "

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "0.9997"
$ws.Cells.Item(22, 5).Value = "  -0.06%  "

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "0.9999"
$ws.Cells.Item(23, 5).Value = "  +0.00%  "

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "6.950"
$ws.Cells.Item(24, 5).Value = "  +11.51%  "

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "9.681"
$ws.Cells.Item(25, 5).Value = "  +6.97%  "

# Row 26
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "0.1548"
$ws.Cells.Item(26, 5).Value = "  +61.78%  "

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "165.42"
$ws.Cells.Item(27, 5).Value = "  +1.30%  "

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "19.84"
$ws.Cells.Item(28, 5).Value = "  +6.26%  "

# Row 29
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "2.203"
$ws.Cells.Item(29, 5).Value = "  +17.35%  "

# Row 30
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "1.571"
$ws.Cells.Item(30, 5).Value = "  +6.97%  "

# Row 31
$ws.Cells.Item(31, 2).Value = "Toncoin"
$ws.Cells.Item(31, 3).Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "1.355"
$ws.Cells.Item(31, 5).Value = "  +2.97%  "

# Row 32
$ws.Cells.Item(32, 2).Value = "Filecoin"
$ws.Cells.Item(32, 3).Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "4.583"
$ws.Cells.Item(32, 5).Value = "  +8.70%  "

# Row 33
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "4.334"
$ws.Cells.Item(33, 5).Value = "  +5.93%  "

# Row 34
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "0.05204"
$ws.Cells.Item(34, 5).Value = "  +8.31%  "

# Row 35
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "1.221"
$ws.Cells.Item(35, 5).Value = "  +9.14%  "

# Row 36
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "0.7516"
$ws.Cells.Item(36, 5).Value = "  +9.34%  "

# Row 37
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "2.774"
$ws.Cells.Item(37, 5).Value = "  +2.50%  "

# Row 38
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "0.9990"
$ws.Cells.Item(38, 5).Value = "  -0.12%  "

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "0.02013"
$ws.Cells.Item(39, 5).Value = "  +7.70%  "

# Row 40
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "2.935"
$ws.Cells.Item(40, 5).Value = "  +4.39%  "

# Row 41
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "6.654"
$ws.Cells.Item(41, 5).Value = "  +6.84%  "

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "78.77"
$ws.Cells.Item(42, 5).Value = "  +6.08%  "

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "0.4682"
$ws.Cells.Item(43, 5).Value = "  +10.55%  "

# Row 44
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "2.080"
$ws.Cells.Item(44, 5).Value = "  +7.52%  "

# Row 45
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "107.23"
$ws.Cells.Item(45, 5).Value = "  +6.40%  "

# Row 46
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "0.8575"
$ws.Cells.Item(46, 5).Value = "  +4.05%  "

# Row 47
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "1.0000"
$ws.Cells.Item(47, 5).Value = "  +0.07%  "

# Row 48
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "9.978"
$ws.Cells.Item(48, 5).Value = "  +4.38%  "

# Row 49
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "7.517"
$ws.Cells.Item(49, 5).Value = "  +8.69%  "

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "0.4290"
$ws.Cells.Item(50, 5).Value = "  +9.86%  "

# Row 51
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "36.55"
$ws.Cells.Item(51, 5).Value = "  +4.12%  "
